$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Docente(s) Responsável(eis)")) {
        $target = $p
        break
    }
}

# Insert a brand-new paragraph right after that heading (inherits the
# heading's formatting for now; we restyle it below).
$target.Range.InsertParagraphAfter() | Out-Null
$inserted = $target.Next()
$inserted.Style = "ListBullet"

# Type both names into the new paragraph, separated by a manual line break
# (Shift+Enter, i.e. a vertical-tab char -> <w:br/> when serialized), the
# same convention already used elsewhere in this document (e.g. the
# "Créditos" bullet list).
$r = $inserted.Range
$r.InsertAfter("5817692 - Katia Cristiane Gandolpho Candioto" + [char]11 + "1176388 - Luiz Tadeu Fernandes Eleno")

# Split the text following the line break into its own run (mirroring how
# Word itself stores each manually-broken line as a separate run) by
# nudging a character-formatting property on it and then reverting it.
$fullText = $inserted.Range.Text
$breakIndex = $fullText.IndexOf([char]11)
$secondStart = $inserted.Range.Start + $breakIndex + 1
$secondEnd = $inserted.Range.End - 1
$secondRun = $d.Range($secondStart, $secondEnd)
$secondRun.Bold = $true
$secondRun2 = $d.Range($secondStart, $secondEnd)
$secondRun2.Bold = $false
